# Rename "Data" sheet to "FCM_Data"
$wb = $excel.ActiveWorkbook
$fcm = $wb.Worksheets.Item("Data")
$fcm.Name = "FCM_Data"

# Add new "DAPI_Data" sheet positioned right after FCM_Data
$ws = $wb.Worksheets.Add($null, $fcm)
$ws.Name = "DAPI_Data"

# Populate header row
$ws.Range("A1").Value = "Treatment"
$ws.Range("B1").Value = "Timepoint"
$ws.Range("C1").Value = "Cells_mL"
$ws.Range("D1").Value = "Cells_mL_Stdev"
$ws.Range("E1").Value = "Mean_Biovolume_um3_cell"
$ws.Range("F1").Value = "Biovolume_Stdev_um3_cell"

# Populate data rows
    $ws.Range("A2").Value = "Control"
    $ws.Range("B2").Value = 0
    $ws.Range("C2").Value = 660667.0
    $ws.Range("D2").Value = 73217.75574473696178756654
    $ws.Range("E2").Value = 0.0455620925016649303
    $ws.Range("F2").Value = 0.00605480485186889444
    $ws.Range("A3").Value = "Control"
    $ws.Range("B3").Value = 4
    $ws.Range("C3").Value = 919405.58319999999366700649
    $ws.Range("D3").Value = 363326.27392113575479015708
    $ws.Range("E3").Value = 0.05080353179310870271
    $ws.Range("F3").Value = 0.01100036917471313574
    $ws.Range("A4").Value = "Control"
    $ws.Range("B4").Value = 8
    $ws.Range("C4").Value = 1133869.72399999992921948433
    $ws.Range("D4").Value = 99930.05014348248369060457
    $ws.Range("E4").Value = 0.04093212299220598538
    $ws.Range("F4").Value = 0.00468449505458882597
    $ws.Range("A5").Value = "Kelp Exudate"
    $ws.Range("B5").Value = 0
    $ws.Range("C5").Value = 663088.10880000004544854164
    $ws.Range("D5").Value = 113546.26705004731775261462
    $ws.Range("E5").Value = 0.03871498619732444135
    $ws.Range("F5").Value = 0.00544648526503371463
    $ws.Range("A6").Value = "Kelp Exudate"
    $ws.Range("B6").Value = 4
    $ws.Range("C6").Value = 1043597.82800000032875686884
    $ws.Range("D6").Value = 181810.62371728930156677961
    $ws.Range("E6").Value = 0.06811629158518681115
    $ws.Range("F6").Value = 0.01349231671113272635
    $ws.Range("A7").Value = "Kelp Exudate"
    $ws.Range("B7").Value = 8
    $ws.Range("C7").Value = 1115268.24240000033751130104
    $ws.Range("D7").Value = 149497.94850320072146132588
    $ws.Range("E7").Value = 0.03272024908449591707
    $ws.Range("F7").Value = 0.00969713437558491695
    $ws.Range("A8").Value = "Kelp Exudate_Nitrate_Phosphate"
    $ws.Range("B8").Value = 0
    $ws.Range("C8").Value = 672115.29839999997057020664
    $ws.Range("D8").Value = 71870.4841418611176777631
    $ws.Range("E8").Value = 0.03630970627897990577
    $ws.Range("F8").Value = 0.0050124924763695301
    $ws.Range("A9").Value = "Kelp Exudate_Nitrate_Phosphate"
    $ws.Range("B9").Value = 4
    $ws.Range("C9").Value = 1226603.58080000011250376701
    $ws.Range("D9").Value = 153931.80801620887359604239
    $ws.Range("E9").Value = 0.10237022179025781488
    $ws.Range("F9").Value = 0.01478012270574448686
    $ws.Range("A10").Value = "Kelp Exudate_Nitrate_Phosphate"
    $ws.Range("B10").Value = 8
    $ws.Range("C10").Value = 1391554.95439999969676136971
    $ws.Range("D10").Value = 81241.07329684446449391544
    $ws.Range("E10").Value = 0.06133378995467635153
    $ws.Range("F10").Value = 0.00992733208000250514
    $ws.Range("A11").Value = "Glucose_Nitrate_Phosphate"
    $ws.Range("B11").Value = 0
    $ws.Range("C11").Value = 646948.58800000010523945093
    $ws.Range("D11").Value = 126328.47284772344573866576
    $ws.Range("E11").Value = 0.05185419746321871109
    $ws.Range("F11").Value = 0.00994342489947468676
    $ws.Range("A12").Value = "Glucose_Nitrate_Phosphate"
    $ws.Range("B12").Value = 4
    $ws.Range("C12").Value = 1665926.80799999949522316456
    $ws.Range("D12").Value = 184066.73734873536159284413
    $ws.Range("E12").Value = 0.09532373834231459842
    $ws.Range("F12").Value = 0.01323736710127628065
    $ws.Range("A13").Value = "Glucose_Nitrate_Phosphate"
    $ws.Range("B13").Value = 8
    $ws.Range("C13").Value = 1544196.52399999974295496941
    $ws.Range("D13").Value = 116715.17866666665941011161
    $ws.Range("E13").Value = 0.04510939859185660039
    $ws.Range("F13").Value = 0.00975322583574211172

# Set selection/scroll position on the Metadata sheet (previously active tab,
# scrolled further down working in column G)
$meta = $wb.Worksheets.Item("Metadata")
$meta.Select()
$excel.ActiveWindow.ScrollRow = 39
$meta.Range("G62").Select()

# Set selection on the new DAPI_Data sheet and make it the active tab
$ws.Select()
$ws.Range("A6").Select()

Write-Output "done"
